$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (shifts existing rows 25.. down by one)
$ws.Rows("25:25").Insert()

# Populate the newly inserted row 25 with the new data record
$ws.Range("A25").Value = 4
$ws.Range("B25").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C25").Value = "Los Lagos"
$ws.Range("D25").Value = 45281
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 300000000
$ws.Range("G25").Value = "Espárragos"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 2000
$ws.Range("N25").Value = "$/kilo"
$ws.Range("O25").Value = "Provincia de Linares"
$ws.Range("P25").Value = 2000
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = "Hortaliza"
